$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: copy formatting from row 4 (B:E only, since A10 has no override)
$ws.Range("B4:E4").Copy()
$ws.Range("B10:E10").PasteSpecial(-4122)

# Row 11: copy formatting from row 5 (A:E, A5 is an empty styled cell)
$ws.Range("A5:E5").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)

# Row 12: copy formatting from row 6 (A:E, A6 is an empty styled cell)
$ws.Range("A6:E6").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)

# Row heights matching the other 21.6pt rows
$ws.Rows.Item(10).RowHeight = 21.6
$ws.Rows.Item(11).RowHeight = 21.6
$ws.Rows.Item(12).RowHeight = 21.6

# Numeric "line number" column values
$ws.Range("B10").Value = 243
$ws.Range("B11").Value = 246
$ws.Range("B12").Value = 205

# New translated text content - added in the same order as the original
# author typed them (column by column for the paired rows 10-11, then
# row by row for row 12), so the shared-string table indices line up.
$ws.Range("C10").Value = ' Ah, welcome back, [hero]\nand [partner]!'
$ws.Range("C11").Value = ' Ahhh! We are happy to see you!\nWe were so very worried about you!'
$ws.Range("D10").Value = ' Ах, с возвращением, [hero]\nи [partner]!'
$ws.Range("D11").Value = ' Аххх! Мы так рады вас видеть!\nМы очень за вас переживали!'
$ws.Range("E10").Value = ' Àö, ò âïèâñàþåîéåí, [hero]\né [partner]!'
$ws.Range("E11").Value = ' Àööö! Íú óàë ñàäú âàò âéäåóû!\nÍú ïœåîû èà âàò ðåñåçéâàìé!'
$ws.Range("C12").Value = ' Our thoughts go with you!\nPlease be strong!'
$ws.Range("D12").Value = ' Все наши мысли только о вас!\nБудьте сильны!'
$ws.Range("E12").Value = ' Âòå îàšé íúòìé óïìûëï ï âàò!\nÁôäûóå òéìûîú!'

# Update view state: scrolled down a bit, with E12 as the active selection
$ws.Range("E12").Select()

